$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cells with default/unedited style, used to restore style after
# forcing numeric-looking text values to be stored as literal text.
$refStyleD = $ws.Range("D38").Style
$refStyleE = $ws.Range("E38").Style

# --- Plain text cells (coin names / links) ---
$ws.Range("B9").Value = 'BTSEToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("B10").Value = 'MXToken'
$ws.Range("C10").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("B11").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C11").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("B12").Value = 'WazirX'
$ws.Range("C12").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("B13").Value = 'MandalaExchangeToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("B14").Value = 'BitrueCoin'
$ws.Range("C14").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("B15").Value = 'BitMartToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("B16").Value = 'BitForexToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("B17").Value = 'TigerCash'
$ws.Range("C17").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("B18").Value = 'LEO'
$ws.Range("C18").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'

# --- Numeric-looking cells (prices / percentages) kept as literal text ---
# Leading apostrophe forces Excel to store the literal text instead of
# re-interpreting it as a number/percentage; the style is then restored
# to the sheet default so no stray number-format/quote-prefix is left behind.
$ws.Range("D2").Value = "'311.28"
$ws.Range("D2").Style = $refStyleD
$ws.Range("E2").Value = "'-1.43%"
$ws.Range("E2").Style = $refStyleE
$ws.Range("D3").Value = "'38.13"
$ws.Range("D3").Style = $refStyleD
$ws.Range("E3").Value = "'-3.12%"
$ws.Range("E3").Style = $refStyleE
$ws.Range("D4").Value = "'5.069"
$ws.Range("D4").Style = $refStyleD
$ws.Range("E4").Value = "'-1.25%"
$ws.Range("E4").Style = $refStyleE
$ws.Range("D5").Value = "'0.07754"
$ws.Range("D5").Style = $refStyleD
$ws.Range("E5").Value = "'-5.17%"
$ws.Range("E5").Style = $refStyleE
$ws.Range("D6").Value = "'4.356"
$ws.Range("D6").Style = $refStyleD
$ws.Range("D7").Value = "'1.898"
$ws.Range("D7").Style = $refStyleD
$ws.Range("E7").Value = "'-4.61%"
$ws.Range("E7").Style = $refStyleE
$ws.Range("D8").Value = "'8.200"
$ws.Range("D8").Style = $refStyleD
$ws.Range("E8").Value = "'-1.48%"
$ws.Range("E8").Style = $refStyleE
$ws.Range("D9").Value = "'3.087"
$ws.Range("D9").Style = $refStyleD
$ws.Range("E9").Value = "'-6.44%"
$ws.Range("E9").Style = $refStyleE
$ws.Range("D10").Value = "'0.9204"
$ws.Range("D10").Style = $refStyleD
$ws.Range("E10").Value = "'-2.03%"
$ws.Range("E10").Style = $refStyleE
$ws.Range("D11").Value = "'0.1227"
$ws.Range("D11").Style = $refStyleD
$ws.Range("E11").Value = "'-5.59%"
$ws.Range("E11").Style = $refStyleE
$ws.Range("D12").Value = "'0.1885"
$ws.Range("D12").Style = $refStyleD
$ws.Range("E12").Value = "'-4.47%"
$ws.Range("E12").Style = $refStyleE
$ws.Range("D13").Value = "'0.08822"
$ws.Range("D13").Style = $refStyleD
$ws.Range("E13").Value = "'-2.19%"
$ws.Range("E13").Style = $refStyleE
$ws.Range("D14").Value = "'0.03414"
$ws.Range("D14").Style = $refStyleD
$ws.Range("E14").Value = "'-2.36%"
$ws.Range("E14").Style = $refStyleE
$ws.Range("D15").Value = "'0.09698"
$ws.Range("D15").Style = $refStyleD
$ws.Range("E15").Value = "'-0.60%"
$ws.Range("E15").Style = $refStyleE
$ws.Range("D16").Value = "'0.001363"
$ws.Range("D16").Style = $refStyleD
$ws.Range("E16").Value = "'-3.15%"
$ws.Range("E16").Style = $refStyleE
$ws.Range("D17").Value = "'0.006022"
$ws.Range("D17").Style = $refStyleD
$ws.Range("E17").Value = "'0.79%"
$ws.Range("E17").Style = $refStyleE
$ws.Range("D18").Value = "'3.559"
$ws.Range("D18").Style = $refStyleD
$ws.Range("E18").Value = "'-2.09%"
$ws.Range("E18").Style = $refStyleE
$ws.Range("D19").Value = "'0.3410"
$ws.Range("D19").Style = $refStyleD
$ws.Range("D20").Value = "'5.026"
$ws.Range("D20").Style = $refStyleD
$ws.Range("E20").Value = "'1.50%"
$ws.Range("E20").Style = $refStyleE
$ws.Range("E21").Value = "'-2.59%"
$ws.Range("E21").Style = $refStyleE
$ws.Range("D22").Value = "'0.2620"
$ws.Range("D22").Style = $refStyleD
$ws.Range("E22").Value = "'1.55%"
$ws.Range("E22").Style = $refStyleE
$ws.Range("D23").Value = "'0.02105"
$ws.Range("D23").Style = $refStyleD
$ws.Range("E23").Value = "'5,592.33%"
$ws.Range("E23").Style = $refStyleE
$ws.Range("E24").Value = "'0.89%"
$ws.Range("E24").Style = $refStyleE
$ws.Range("E25").Value = "'-2.39%"
$ws.Range("E25").Style = $refStyleE
$ws.Range("D26").Value = "'0.004253"
$ws.Range("D26").Style = $refStyleD
$ws.Range("E26").Value = "'-10.84%"
$ws.Range("E26").Style = $refStyleE
$ws.Range("D27").Value = "'0.0001351"
$ws.Range("D27").Style = $refStyleD
$ws.Range("E27").Value = "'-65.29%"
$ws.Range("E27").Style = $refStyleE
$ws.Range("D39").Value = "'0.02135"
$ws.Range("D39").Style = $refStyleD
$ws.Range("E39").Value = "'-3.37%"
$ws.Range("E39").Style = $refStyleE
$ws.Range("D40").Value = "'0.05008"
$ws.Range("D40").Style = $refStyleD
$ws.Range("E40").Value = "'-3.57%"
$ws.Range("E40").Style = $refStyleE
$ws.Range("D41").Value = "'0.007765"
$ws.Range("D41").Style = $refStyleD
$ws.Range("E41").Value = "'0.17%"
$ws.Range("E41").Style = $refStyleE
$ws.Range("D42").Value = "'0.009958"
$ws.Range("D42").Style = $refStyleD
$ws.Range("E42").Value = "'-3.69%"
$ws.Range("E42").Style = $refStyleE
$ws.Range("E43").Value = "'-4.11%"
$ws.Range("E43").Style = $refStyleE
$ws.Range("D44").Value = "'0.002061"
$ws.Range("D44").Style = $refStyleD
$ws.Range("E44").Value = "'-1.88%"
$ws.Range("E44").Style = $refStyleE
$ws.Range("D45").Value = "'0.009689"
$ws.Range("D45").Style = $refStyleD
$ws.Range("E45").Value = "'4.53%"
$ws.Range("E45").Style = $refStyleE
$ws.Range("D46").Value = "'0.00006466"
$ws.Range("D46").Style = $refStyleD
$ws.Range("E46").Value = "'-4.10%"
$ws.Range("E46").Style = $refStyleE
$ws.Range("E47").Value = "'0.03%"
$ws.Range("E47").Style = $refStyleE
$ws.Range("D48").Value = "'0.003210"
$ws.Range("D48").Style = $refStyleD
$ws.Range("E48").Value = "'11.32%"
$ws.Range("E48").Style = $refStyleE
$ws.Range("E49").Value = "'-0.11%"
$ws.Range("E49").Style = $refStyleE
$ws.Range("E50").Value = "'0.03%"
$ws.Range("E50").Style = $refStyleE
$ws.Range("E51").Value = "'0.03%"
$ws.Range("E51").Style = $refStyleE
